$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.131.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.318.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("E7").Value = "  +4.14%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.315.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.889.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.132"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.178.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.315.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "425.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.38%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.460.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.513"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("E28").Value = "  +5.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.65%  "

$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("E37").Value = "  -4.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.878.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.74%  "

$ws.Range("E43").Value = "  -4.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0661"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.19%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "313.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("E51").Value = "  -0.81%  "
